$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 146 entirely ("「私は〜する」" post) — remaining rows shift up by one.
$ws.Rows.Item(146).Delete()
